$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "Z1" = 0.98399148769620681
    "BC1" = 0.74675569609372539
    "O3" = 0.90835820592323646
    "W3" = 0.89087508360014134
    "BI4" = 0.83331642805799722
    "V5" = 0.86518857106124059
    "AC5" = 0.96245361077412106
    "BD5" = 0.67219452734216434
    "D6" = 0.73108148630698899
    "AO6" = 0.91262221733806759
    "AP6" = 0.97648722024043644
    "AU7" = 0.74239593972364126
    "BL7" = 0.90679331431840393
    "N8" = 0.8105151472441714
    "AV8" = 0.85749525724020426
    "BB8" = 0.89299342847351526
    "BO9" = 0.91774567074488811
    "B10" = 0.68193760169130568
    "V10" = 0.88984175106980712
    "AR10" = 0.81130832467837011
    "BM11" = 0.99286705748433945
    "BC12" = 0.91988991074523208
    "N13" = 0.70716798170895934
    "AH13" = 0.77908416940005121
    "BP13" = 0.62098491905961006
    "X14" = 0.80886375987649783
    "AP14" = 0.67394477927074936
    "BB15" = 0.84620273564782489
    "BN15" = 0.87617365572012229
    "O16" = 0.76956887644675798
    "T16" = 0.55739726621804642
    "BP16" = 0.83255611105081906
    "AP17" = 0.65870628264236641
    "AQ17" = 0.59069387858616085
    "T18" = 0.72675073387988287
    "AS18" = 0.86350575739266056
    "D19" = 0.72578295128163484
    "AN19" = 0.96926962238045999
    "F20" = 0.64186284731711751
    "BC20" = 0.9651967867517155
    "A21" = 0.74304159966320082
    "C21" = 0.67023133125274881
    "V21" = 0.88841280903183639
    "AI22" = 0.95044538945517232
    "AJ22" = 0.73036621409775604
    "O23" = 0.94555717881762846
    "Y23" = 0.87319460056856602
    "I24" = 0.99760084168153296
    "AI24" = 0.83360461650939466
    "AY24" = 0.816899939532999
    "F25" = 0.73876549167125605
    "AZ26" = 0.75848781236005036
    "BO26" = 0.82161195897777972
    "O27" = 0.86097139000049228
    "AP27" = 0.77738030028057581
    "BH27" = 0.69344303574366584
    "AU28" = 0.88766466197910199
    "AI29" = 0.85493868572681941
    "P31" = 0.73643988147019579
    "AL31" = 0.91717143593434713
    "AQ31" = 0.97407921396202179
    "J32" = 0.91484905115717274
    "AD32" = 0.8170245060854342
    "U34" = 0.64721459943499005
    "AD34" = 0.98931644534776664
    "AG34" = 0.77962962846194594
    "AW35" = 0.75277012824104594
    "BF35" = 0.93015286849337553
    "AK36" = 0.87515165278945273
    "BC36" = 0.75512657755861778
    "AB37" = 0.69867806271202304
    "AM37" = 0.9872906883303284
    "AS37" = 0.8988165470377687
    "L39" = 0.60020540967537483
    "AT39" = 0.99658313944664811
    "M40" = 0.81982708618129885
    "U40" = 0.94979387898127543
    "AZ40" = 0.91519281148914611
    "Y41" = 0.72876655166827142
    "AG41" = 0.75217756630398225
    "G42" = 0.68072824971327672
    "BC42" = 0.95177012001883121
    "BP42" = 0.74627076632954947
    "AD43" = 0.69295530988754983
    "AQ44" = 0.87567553914793583
    "AT44" = 0.88120026805077589
    "BK44" = 0.96994712607445832
    "K45" = 0.80658666641608123
    "V45" = 0.58033999006957537
    "AX45" = 0.9759203010927584
    "K46" = 0.99737460580584725
    "V46" = 0.68940972742874984
    "I47" = 0.59803485651367061
    "AF48" = 0.74792565776657294
    "AU48" = 0.73782002885575526
    "S49" = 0.74569876745816788
    "AU49" = 0.98784445175581137
    "G50" = 0.74727371498108441
    "M50" = 0.98677585987791039
    "BF50" = 0.87260667177068685
    "AO51" = 0.99518519922514215
    "A53" = 0.96766385133053134
    "U53" = 0.79435254923501764
    "AW53" = 0.91441754032818978
    "AY53" = 0.81015323371273329
    "AL54" = 0.89700121066433058
    "B56" = 0.83452547541152355
    "AA56" = 0.69665028738726975
    "AY56" = 0.89449223561730773
    "P57" = 0.69260182494445677
    "R57" = 0.94757449158233387
    "Z57" = 0.91960313326891363
    "AM57" = 0.84899312924570336
    "AV57" = 0.7994716340283039
    "BH58" = 0.94908426211442387
    "AI59" = 0.88983344235033457
    "AT59" = 0.78981946068696929
    "U60" = 0.88857298012154673
    "BG60" = 0.99191516818575676
    "BL60" = 0.87345822024126019
    "BH61" = 0.91049734292143636
    "BJ61" = 0.94981919316146624
    "H62" = 0.96321672871520292
    "BP62" = 0.92767187149806918
    "Y63" = 0.99735872979316131
    "AQ63" = 0.86286157678700803
    "BI63" = 0.67851921410867821
    "D64" = 0.98105085874765918
    "AK64" = 0.80737470205576622
    "BJ64" = 0.97739328028959094
    "BM64" = 0.92842131853443166
    "AH65" = 0.94243405141227599
    "BI65" = 0.75420011306488621
    "AO66" = 0.90856192961449822
    "BB66" = 0.82661418415005028
    "M67" = 0.87660885579158043
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
